$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Mark three more milestone rows as "III" / completed ("X"), matching the
# pattern already used on other rows (e.g. row 5: E="I", F="X").
$ws.Range("E20").Value = "III"
$ws.Range("F20").Value = "X"

$ws.Range("E55").Value = "III"
$ws.Range("F55").Value = "X"

$ws.Range("E75").Value = "III"
$ws.Range("F75").Value = "X"

# Reflect the final selection/active cell left on the sheet (F55), with no
# frozen/scrolled topLeftCell override.
$ws.Range("F55").Select()

$wb.Save()
